$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 24,13
$data[0,0] = 1.489221558399606
$data[0,1] = 0.06233418494747411
$data[0,2] = 0.02761335154428224
$data[0,3] = 0
$data[0,4] = 3.449164834304483
$data[0,5] = 0.00260577339768102
$data[0,6] = 0
$data[0,7] = 2.139717580206479
$data[0,8] = 0.1694091139184017
$data[0,9] = 1.261658422223206
$data[0,10] = 0.416834724591709
$data[0,11] = 0
$data[0,12] = 3.405052693263585
$data[1,0] = 1.453414093748847
$data[1,1] = 0.05677461439502451
$data[1,2] = 0.02716868039836839
$data[1,3] = 0
$data[1,4] = 3.447532439062783
$data[1,5] = 0.002609703877810129
$data[1,6] = 0
$data[1,7] = 2.14287825725755
$data[1,8] = 0.1699293051653861
$data[1,9] = 1.222573582026769
$data[1,10] = 0.4129758485621835
$data[1,11] = 0
$data[1,12] = 3.424097129044696
$data[2,0] = 1.432214001741812
$data[2,1] = 0.05339267125553704
$data[2,2] = 0.02690628017263563
$data[2,3] = 0
$data[2,4] = 3.447893372780314
$data[2,5] = 0.002612246616642539
$data[2,6] = 0
$data[2,7] = 2.145598864110951
$data[2,8] = 0.1702953352610734
$data[2,9] = 1.199276176095651
$data[2,10] = 0.4107949931515691
$data[2,11] = 0
$data[2,12] = 3.436615970012383
$data[3,0] = 1.423772773556124
$data[3,1] = 0.05202239869191772
$data[3,2] = 0.02680204452699542
$data[3,3] = 0
$data[3,4] = 3.448383384429263
$data[3,5] = 0.002613315445956853
$data[3,6] = 0
$data[3,7] = 2.146903718303079
$data[3,8] = 0.1704562414844943
$data[3,9] = 1.189958732117503
$data[3,10] = 0.409953763695782
$data[3,11] = 0
$data[3,12] = 3.441925057096249
$data[4,0] = 1.422383082853457
$data[4,1] = 0.05179534108093264
$data[4,2] = 0.02678489966505282
$data[4,3] = 0
$data[4,4] = 3.448485466541356
$data[4,5] = 0.002613494898590904
$data[4,6] = 0
$data[4,7] = 2.147132239345147
$data[4,8] = 0.1704836698699097
$data[4,9] = 1.18842224025596
$data[4,10] = 0.4098169487783565
$data[4,11] = 0
$data[4,12] = 3.442819162533027
$data[5,0] = 1.432099358266981
$data[5,1] = 0.05337415938971901
$data[5,2] = 0.02690486347418286
$data[5,3] = 0
$data[5,4] = 3.447898592577886
$data[5,5] = 0.002612260898899986
$data[5,6] = 0
$data[5,7] = 2.145615667409011
$data[5,8] = 0.1702974577123193
$data[5,9] = 1.199149803122907
$data[5,10] = 0.4107834556515542
$data[5,11] = 0
$data[5,12] = 3.436686729857705
$data[6,0] = 1.47671225391457
$data[6,1] = 0.06041063492597232
$data[6,2] = 0.02745783549200453
$data[6,3] = 0
$data[6,4] = 3.448319116344052
$data[6,5] = 0.002607101830158084
$data[6,6] = 0
$data[6,7] = 2.140645539247551
$data[6,8] = 0.1695788097598658
$data[6,9] = 1.248036615503025
$data[6,10] = 0.4154651071276021
$data[6,11] = 0
$data[6,12] = 3.411447876894592
$data[7,0] = 1.570423254052571
$data[7,1] = 0.07446391638528382
$data[7,2] = 0.02862568696967926
$data[7,3] = 0
$data[7,4] = 3.459958352958864
$data[7,5] = 0.002598007051293458
$data[7,6] = 0
$data[7,7] = 2.137086491275447
$data[7,8] = 0.1685386934743498
$data[7,9] = 1.349462152124147
$data[7,10] = 0.4261388077752599
$data[7,11] = 0
$data[7,12] = 3.368504960768149
$data[8,0] = 1.643064025855324
$data[8,1] = 0.08495023558130299
$data[8,2] = 0.02953355680278946
$data[8,3] = 0
$data[8,4] = 3.475105468375645
$data[8,5] = 0.002591941722973137
$data[8,6] = 0
$data[8,7] = 2.138244631171688
$data[8,8] = 0.1679985078660096
$data[8,9] = 1.427375231913459
$data[8,10] = 0.4348886806815528
$data[8,11] = 0
$data[8,12] = 3.340948712825067
$data[9,0] = 1.676933330408701
$data[9,1] = 0.08975714578832594
$data[9,2] = 0.02995718497299293
$data[9,3] = 0
$data[9,4] = 3.483429592252776
$data[9,5] = 0.002589314950370005
$data[9,6] = 0
$data[9,7] = 2.139591206004454
$data[9,8] = 0.1678011847547189
$data[9,9] = 1.463559469644196
$data[9,10] = 0.439065993177735
$data[9,11] = 0
$data[9,12] = 3.329280125310675
$data[10,0] = 1.689877129948172
$data[10,1] = 0.09158275276573136
$data[10,2] = 0.03011911143151735
$data[10,3] = 0
$data[10,4] = 3.486787874414858
$data[10,5] = 0.002588339190346833
$data[10,6] = 0
$data[10,7] = 2.140218997039327
$data[10,8] = 0.1677334065037748
$data[10,9] = 1.477368053965336
$data[10,10] = 0.4406760952134334
$data[10,11] = 0
$data[10,12] = 3.324986242883028
$data[11,0] = 1.687084198464788
$data[11,1] = 0.09118933731946299
$data[11,2] = 0.03008417095742288
$data[11,3] = 0
$data[11,4] = 3.486055441672008
$data[11,5] = 0.00258854849679313
$data[11,6] = 0
$data[11,7] = 2.140078548631294
$data[11,8] = 0.1677476952167005
$data[11,9] = 1.47438939901852
$data[11,10] = 0.4403280759682673
$data[11,11] = 0
$data[11,12] = 3.325905456969195
$data[12,0] = 1.677995857120663
$data[12,1] = 0.08990723227933017
$data[12,2] = 0.02997047666686825
$data[12,3] = 0
$data[12,4] = 3.483701749444293
$data[12,5] = 0.002589234295209822
$data[12,6] = 0
$data[12,7] = 2.139640492425372
$data[12,8] = 0.1677954695345321
$data[12,9] = 1.464693378830049
$data[12,10] = 0.4391978916181785
$data[12,11] = 0
$data[12,12] = 3.328924364081416
$data[13,0] = 1.672444372063921
$data[13,1] = 0.08912260261183746
$data[13,2] = 0.02990103134060007
$data[13,3] = 0
$data[13,4] = 3.482286886530801
$data[13,5] = 0.002589656830324157
$data[13,6] = 0
$data[13,7] = 2.139387520744833
$data[13,8] = 0.1678256364433786
$data[13,9] = 1.45876813740179
$data[13,10] = 0.4385092970629074
$data[13,11] = 0
$data[13,12] = 3.33078978215508
$data[14,0] = 1.640867150817371
$data[14,1] = 0.08463683551437384
$data[14,2] = 0.02950608384290376
$data[14,3] = 0
$data[14,4] = 3.474590312272227
$data[14,5] = 0.00259211604450491
$data[14,6] = 0
$data[14,7] = 2.138173124017982
$data[14,8] = 0.1680123760230074
$data[14,9] = 1.425025405719538
$data[14,10] = 0.4346196400568374
$data[14,11] = 0
$data[14,12] = 3.341728734996622
$data[15,0] = 1.621706507714862
$data[15,1] = 0.08189439224948103
$data[15,2] = 0.02926650455250268
$data[15,3] = 0
$data[15,4] = 3.470235871492378
$data[15,5] = 0.002593658530093714
$data[15,6] = 0
$data[15,7] = 2.137638080018348
$data[15,8] = 0.1681393211181152
$data[15,9] = 1.404515037010214
$data[15,10] = 0.4322838504608058
$data[15,11] = 0
$data[15,12] = 3.348661529562136
$data[16,0] = 1.610763453781658
$data[16,1] = 0.08032045527875198
$data[16,2] = 0.02912970703788176
$data[16,3] = 0
$data[16,4] = 3.46786625372097
$data[16,5] = 0.002594558192488039
$data[16,6] = 0
$data[16,7] = 2.137407493862241
$data[16,8] = 0.1682168942006719
$data[16,9] = 1.392787790108969
$data[16,10] = 0.4309589086146275
$data[16,11] = 0
$data[16,12] = 3.352730678862656
$data[17,0] = 1.607071669596394
$data[17,1] = 0.07978813671907403
$data[17,2] = 0.02908356252758892
$data[17,3] = 0
$data[17,4] = 3.467087119383137
$data[17,5] = 0.002594864946519078
$data[17,6] = 0
$data[17,7] = 2.137342673321172
$data[17,8] = 0.1682439424273561
$data[17,9] = 1.388829135327114
$data[17,10] = 0.430513493128899
$data[17,11] = 0
$data[17,12] = 3.354122432937302
$data[18,0] = 1.623738157895218
$data[18,1] = 0.08218597310728626
$data[18,2] = 0.02929190461395592
$data[18,3] = 0
$data[18,4] = 3.470685443488051
$data[18,5] = 0.002593493040466105
$data[18,6] = 0
$data[18,7] = 2.137687050935668
$data[18,8] = 0.1681253360271739
$data[18,9] = 1.406691180692263
$data[18,10] = 0.4325305807276294
$data[18,11] = 0
$data[18,12] = 3.347915076833232
$data[19,0] = 1.680662117386134
$data[19,1] = 0.09028367229910828
$data[19,2] = 0.03000383068829393
$data[19,3] = 0
$data[19,4] = 3.484387492211013
$data[19,5] = 0.002589032346216214
$data[19,6] = 0
$data[19,7] = 2.139765961145883
$data[19,8] = 0.1677812487472643
$data[19,9] = 1.467538448799274
$data[19,10] = 0.439529088220425
$data[19,11] = 0
$data[19,12] = 3.328034251148452
$data[20,0] = 1.71855418284332
$data[20,1] = 0.09560711117336496
$data[20,2] = 0.03047789355485975
$data[20,3] = 0
$data[20,4] = 3.494543968189831
$data[20,5] = 0.002586227387229567
$data[20,6] = 0
$data[20,7] = 2.14181169982929
$data[20,8] = 0.1675968331543629
$data[20,9] = 1.507925747626103
$data[20,10] = 0.4442676126160023
$data[20,11] = 0
$data[20,12] = 3.315768219785127
$data[21,0] = 1.698267550012361
$data[21,1] = 0.09276302156104066
$data[21,2] = 0.03022408112012442
$data[21,3] = 0
$data[21,4] = 3.489013348789641
$data[21,5] = 0.002587714379204411
$data[21,6] = 0
$data[21,7] = 2.140656985953484
$data[21,8] = 0.167691562408109
$data[21,9] = 1.486313607099817
$data[21,10] = 0.4417235376961202
$data[21,11] = 0
$data[21,12] = 3.322248254827002
$data[22,0] = 1.622819422505415
$data[22,1] = 0.0820541409042761
$data[22,2] = 0.02928041831768269
$data[22,3] = 0
$data[22,4] = 3.47048177516352
$data[22,5] = 0.00259356781831901
$data[22,6] = 0
$data[22,7] = 2.137664671264588
$data[22,8] = 0.1681316443860581
$data[22,9] = 1.405707145403113
$data[22,10] = 0.4324189781048773
$data[22,11] = 0
$data[22,12] = 3.348252288328737
$data[23,0] = 1.544405897616286
$data[23,1] = 0.07063423156432691
$data[23,2] = 0.02830092385741523
$data[23,3] = 0
$data[23,4] = 3.455651574272721
$data[23,5] = 0.002600358679022489
$data[23,6] = 0
$data[23,7] = 2.137386817539401
$data[23,8] = 0.1687806622331856
$data[23,9] = 1.321428054108196
$data[23,10] = 0.4230916907515763
$data[23,11] = 0
$data[23,12] = 3.379420845064047

$ws.Range("B2:N25").Value = $data
